$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Strip the "Q11.x. " / "Q14. " survey-question prefixes from the skill
# names in column A (and normalize "Computer Skills" casing), now that
# the sheet is being reused for the 2nd-gen tool.
$ws.Range("A2").Value = "Generative AI skills"
$ws.Range("A10").Value = "General computer skills"
$ws.Range("A3").Value = "Market Analysis and Customer Understanding"
$ws.Range("A4").Value = "Creating and Testing Business Ideas"
$ws.Range("A5").Value = "Making a Business Plan"
$ws.Range("A6").Value = "Running a Business"
$ws.Range("A7").Value = "Branding and Marketing"
$ws.Range("A8").Value = "Sales and Customer Service"
$ws.Range("A9").Value = "Future-Thinking"

# Move the selection to B7, matching the saved cursor position.
$ws.Range("B7").Select()
